# Generate Report for Handoff
#
# A fresh handoff XLIFF was generated for the source file
# "7f75ee81-dcb1-4aa8-b86e-028f717cde7e.md" (row 5 on every status sheet).
# Update the recorded handoff timestamps accordingly:
#   - Overview!G5            "Latest HO Xliff Generate Date"
#   - zh-cn!H5                "Latest Handoff Datetime"
#   - de-de!H5                "Latest Handoff Datetime"

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G5").Value = "2016-08-26 18:43:49"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H5").Value = "2016-08-26 18:43:45"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H5").Value = "2016-08-26 18:43:49"
